$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Shorten the "Use noexcept..." bullet to "Use constexpr where appropriate."
# ---------------------------------------------------------------------------
$oldNoexceptText = "Use noexcept where appropriate, and check existing cases where noexcept is being used on functions returning STL objects by value etc. May require consultation of the standard to check which STL container operations are noexcept."
$d.Content.Find.Execute($oldNoexceptText, $false, $false, $false, $false, $false, $true, 1, $false, "Use constexpr where appropriate.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert a brand-new bullet, carrying the original (long) "Use noexcept..."
#    text, right before the "Explicitly qualify calls..." bullet.
# ---------------------------------------------------------------------------
$explicitRng = $d.Content
$explicitRng.Find.Execute("Explicitly qualify calls to functions in the", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$explicitRng.InsertParagraphBefore() | Out-Null

# The paragraph we just inserted sits right before the "Explicitly qualify..." one.
$explicitPara = $explicitRng.Paragraphs(1).Next()
$newPara = $explicitPara.Previous()
$newParaRange = $newPara.Range
# Keep the trailing paragraph mark; only fill in the text that precedes it.
$insertRange = $d.Range($newParaRange.Start, $newParaRange.End - 1)
$insertRange.Text = $oldNoexceptText

# ---------------------------------------------------------------------------
# 3) / 4) Move the two <w:lastRenderedPageBreak/> markers one bullet earlier
#    (repagination side effect): they now land on "Investigate nested
#    exceptions..." and ".NET injection." instead of the bullets that used
#    to follow them.
# ---------------------------------------------------------------------------
function Move-LastRenderedPageBreakToRunStart($doc, $searchText) {
    $rng = $doc.Content
    $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $full = $doc.Range($rng.Start, $rng.End)
    $escaped = $searchText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $full.InsertXML($xmlFrag)
}

function Remove-LastRenderedPageBreakFromRun($doc, $searchText) {
    $rng = $doc.Content
    $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $full = $doc.Range($rng.Start, $rng.End)
    $escaped = $searchText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $full.InsertXML($xmlFrag)
}

# "Reduce compile time." loses its page break...
Remove-LastRenderedPageBreakFromRun $d "Reduce compile time."
# ...which now belongs on "Investigate nested exceptions (std::throw_with_nested etc.)."
Move-LastRenderedPageBreakToRunStart $d "Investigate nested exceptions (std::throw_with_nested etc.)."

# "Without DLL dependency if possible." loses its page break...
Remove-LastRenderedPageBreakFromRun $d "Without DLL dependency if possible."
# ...which now belongs on ".NET injection."
Move-LastRenderedPageBreakToRunStart $d ".NET injection."

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
